$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1, matching the existing header formatting (bold, centered, bordered)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Reuse the exact style already applied to the other header cells (A1:E1)
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean outlier-flag values for rows 2-12 across columns F (KNN), G (SVM), H (RF)
$values = @(
    @(0, 0, 0),  # row 2
    @(1, 0, 1),  # row 3
    @(1, 0, 0),  # row 4
    @(0, 0, 0),  # row 5
    @(0, 0, 0),  # row 6
    @(0, 0, 0),  # row 7
    @(0, 0, 0),  # row 8
    @(0, 0, 0),  # row 9
    @(0, 0, 0),  # row 10
    @(0, 0, 0),  # row 11
    @(0, 0, 0)   # row 12
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $triple = $values[$i]
    $ws.Cells.Item($row, 6).Value = [bool]$triple[0]
    $ws.Cells.Item($row, 7).Value = [bool]$triple[1]
    $ws.Cells.Item($row, 8).Value = [bool]$triple[2]
}
